$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new profit row for 2025-09-14 right after the last existing
# data row (row 27), matching the existing column layout: Date (text) / Profit (number).
$row = 28

# Leading apostrophe forces Excel to store the date as literal text,
# matching the existing date cells in column A (e.g. A27 = "09/13/2025"),
# instead of letting it auto-convert to a date serial value.
$ws.Cells.Item($row, 1).Value = "'09/14/2025"
# Reset to the default/Normal style so the cell doesn't pick up the
# quote-prefix style Excel assigns to forced-text entries, keeping it
# consistent with the unstyled data cells above it.
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 16222.95
